# Update the "Rules" worksheet: change cell E8's text from "Good Morning"
# to "GIT UPDATE", and leave the active selection on E8 (as it was left
# selected in Excel after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
